# Balance iteration 2 - Added second testmap
# Adds a new "Deck-Names" worksheet (mapping Building-List "Type" values to
# their corresponding Deck card names) at the end of the workbook, and makes
# it the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Add the new worksheet after the last existing sheet.
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "Deck-Names"

# Building "Type" (column A) -> Deck card name (column B)
$rows = @(
    @("City",       "City"),
    @("Power",      "Electricity"),
    @("Tavern",     "GettingDrunk"),
    @("Park",       "Park"),
    @("Brewery",    "Booze"),
    @("Tower",      "Outlook"),
    @("Circus",     "Entertainment"),
    @("Brothel",    "Pleasure"),
    @("University", "Knowledge"),
    @("Mine",       "Resource"),
    @("Sawmill",    "Wood"),
    @("Forge",      "Metall"),
    @("Storage",    "Storehouse"),
    @("Temple",     "Belief "),
    @("Fountain",   "Artwork")
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $i + 1
    $ws.Cells.Item($r, 1).Value = $rows[$i][0]
    $ws.Cells.Item($r, 2).Value = $rows[$i][1]
}

# Match the source column-B width (autofit-ish).
$ws.Columns.Item(2).ColumnWidth = 13.86

# Match this workbook's standard (metric / 2cm) top+bottom page margins,
# same as every other sheet in the file.
$ws.PageSetup.TopMargin = 56.692913399999995
$ws.PageSetup.BottomMargin = 56.692913399999995

# Select the whole table, as in the source workbook, and make this the
# active sheet/tab shown when the workbook is opened.
$ws.Range("A1:B15").Select()
$ws.Activate()
